# Soutenance Projet 11 - update presentation date, title and spelling fixes
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Header/footer "datetimeFigureOut" field: 06/05/2019 -> 13/05/2019
#    This cached field text lives on the slide master and on every
#    slide layout (each has its own "Date Placeholder ..." shape).
# ---------------------------------------------------------------------
$newDate = "13/05/2019"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $cl = $layouts.Item($L)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 - title: merge the two runs "Soutenance Projet " + "11"
#    into a single run "Soutenance Projet 11".
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
# Write a throwaway value first so the engine doesn't try to diff the
# old two-run text against the new text (which would keep the split).
$titleRange.Text = "-"
$titleShape.TextFrame.TextRange.Text = "Soutenance Projet 11"

# ---------------------------------------------------------------------
# 3) Slide 1 - subtitle: second paragraph "xx/05/2019" -> "22/05/2019"
# ---------------------------------------------------------------------
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange
$dateParaThrow = $subtitleRange.Paragraphs(2, 1)
$dateParaThrow.Text = "-"
$dateParaFinal = $subtitleShape.TextFrame.TextRange.Paragraphs(2, 1)
$dateParaFinal.Text = "22/05/2019"

# ---------------------------------------------------------------------
# 4) Slide 2 - spelling/content unchanged; the Issue #4 line already
#    reads correctly, the diff only drops the redundant trailing
#    endParaRPr on that paragraph.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$bodyShape = $slide2.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$issue4Throw = $bodyRange.Paragraphs(7, 1)
$issue4Throw.Text = "-"
$issue4Final = $bodyShape.TextFrame.TextRange.Paragraphs(7, 1)
$issue4Final.Text = "Issue #4 : R" + [char]0x00E9 + "initialisation du mot de passe"
